$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Issue number becomes a real number (1) instead of the text "#1"
$ws.Range("A2").Value = 1

# Row 3: new bug - Heatmap / "semanticplacein json output  not sorted"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Heatmap"

# Row 4: new bug - Jsons / "Jsons ordering unable to suit to wiki's requirements"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Jsons"

$ws.Range("C3").Value = "L"
$ws.Range("C4").Value = "L"

# The Raw Score / Final Score columns use a shared formula that already
# spans D2:D4 / E2:E4 - fill it down into the newly used rows, picking up
# the same number format/style as the master (row 2) cells.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D3").Formula = "=VLOOKUP(C3,`$I`$2:`$J`$5,2,FALSE)"
$ws.Range("E3").Formula = '=IF(G3 = "Y",0,D3)'
$ws.Range("D4").Formula = "=VLOOKUP(C4,`$I`$2:`$J`$5,2,FALSE)"
$ws.Range("E4").Formula = '=IF(G4 = "Y",0,D4)'

$ws.Range("F4").Value = "Jsons ordering unable to suit to wiki's requirements"
$ws.Range("F3").Value = "semanticplacein json output  not sorted"

$ws.Range("G3").Value = "Y"
$ws.Range("G4").Value = "Y"

# Widen the Description column so the new text is readable
$ws.Columns.Item(6).ColumnWidth = 46

# Move the active selection to F4 (where the edit ended)
[void]$ws.Range("F4").Select()
